$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary table in columns E:F (rows 4-9), filled before the header row ---
$ws.Range("E4").Value = "Minor"
$ws.Range("F4").Formula = "=COUNTIF(B2:B21, ""<18"")"

$ws.Range("E5").Value = "Major"
$ws.Range("F5").Formula = "=COUNTIFS(B2:B21, "">=18"", B2:B21, ""<=25"")"

$ws.Range("E6").Value = "Middle age"
$ws.Range("F6").Formula = "=COUNTIFS(B2:B21, "">=26"", B2:B21, ""<=40"")"

$ws.Range("E7").Value = "Above middle age"
$ws.Range("F7").Formula = "=COUNTIFS(B2:B21, "">=40"", B2:B21, ""<=60"")"

$ws.Range("E8").Value = "Senior Citizen"
$ws.Range("F8").Formula = "=COUNTIFS(B2:B21, "">=61"", B2:B21, ""<=100"")"

$ws.Range("E9").Value = "Total"
$ws.Range("F9").Formula = "=SUM(F4:F8)"

# --- Header row for the summary table, bold, added last ---
$ws.Range("E3:F3").Font.Bold = $true
$ws.Range("E3").Value = "Category"
$ws.Range("F3").Value = "Counts"

# --- Column E width ---
$ws.Range("E1").ColumnWidth = 15.1666667

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("B3").Select() | Out-Null

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Calculate()
